# "Drop in files from RMI script"
#
# Updates the SoCDTtiNTY-psgr assumptions that were recalibrated by the RMI
# script, and leaves the workbook with the "About" sheet as the active tab
# (matching the view state captured when the script re-saved the file).

$wb = $excel.ActiveWorkbook

$psgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# B2 previously held a helper formula (=0.076+(0.076-0.0725)) that only
# existed to nudge the value to 0.0795; the script replaces it with the
# plain calibrated value used everywhere else in the row.
$psgr.Range("B2").Value = 0.076

# D2 (aircraft) is recalibrated to match the rest of row 2.
$psgr.Range("D2").Value = 0.076

# B5 and E5 (rail, motorbikes placeholder columns) are recalibrated from the
# old 0.01 placeholder to the 0.029 value used across the rest of row 5.
$psgr.Range("B5").Value = 0.029
$psgr.Range("E5").Value = 0.029

# The workbook now opens on the "About" sheet instead of "SoCDTtiNTY-psgr".
$wb.Worksheets.Item("About").Activate()
